$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append new row 3 with the new test-mail entry ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A3").Value = "Kun jij dit even regelen?"
$ws.Range("B3").Value = "mailmind.test@zohomail.eu"
$ws.Range("C3").Value = "Testmail #1: Kun jij dit even regelen?"
$ws.Range("D3").Value = "Overig"
$ws.Range("E3").Value = "Dank voor je bericht. We pakken dit intern op en houden je op de hoogte."
$ws.Range("F3").Value = "2025-07-31 21:23:08"
$ws.Range("G3").Value = "Ja"
$ws.Range("H3").Value = "Ja"
$ws.Range("I3").Value = "Nee"
$ws.Range("J3").Value = "Nee"

# --- Extend the conditional formatting ranges to include the new row ---
$colsToExtend = "D", "G", "H", "I", "J"
foreach ($col in $colsToExtend) {
    $srcCell = $ws.Range($col + "2")
    $newRange = $ws.Range($col + "2:" + $col + "3")
    $rules = $srcCell.FormatConditions
    for ($i = 1; $i -le $rules.Count; $i++) {
        $rules.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Sheet "Dashboard": bump the "Overig" count from 1 to 2 ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 2
